$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete columns C:D (code_insee_commune, code_postal_commune), shifting
# everything to their right left by two columns.
$ws.Range("C1:D1").EntireColumn.Delete()

# Leave the selection where Excel puts it after a column delete: the cell
# that used to be the first deleted column.
[void]$ws.Range("C1").Select()

# Touch the two now-blank trailing cells of row 2 so they stay part of the
# sheet's used range (mirrors Excel re-touching the former L2/M2 cells).
$ws.Range("L2:M2").NumberFormat = "General"
